# Rename the skillTypes value "OFFENSIVE" to "DAMAGE" (and the combined
# "BUFF, OFFENSIVE" to "BUFF, DAMAGE") throughout the skillTypes column (H)
# of the test-skills worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the extent of the used data and the skillTypes column by header text.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

$typesCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value() -eq "skillTypes") {
        $typesCol = $c
        break
    }
}
if ($typesCol -eq 0) { $typesCol = 8 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $typesCol)
    $val = $cell.Value()
    if ($val -eq "OFFENSIVE") {
        $cell.Value = "DAMAGE"
    } elseif ($val -eq "BUFF, OFFENSIVE") {
        $cell.Value = "BUFF, DAMAGE"
    }
}

# Update the active selection to H3, as recorded in the saved sheet view.
$ws.Range("H3").Select()
